# 16.1.3 indicator sheet: the two multilingual section-header rows ("Age" /
# "Education" dividers at rows 18 and 28) get reworded to the "By ..." /
# "По ..." / "... боюнча" phrasing used elsewhere in the workbook.
#
# Row 18 (age-group section divider): "Age (in years)" -> "By age (in years)"
# Row 28 (education section divider): "Education" -> "By education"
#
# The write order below (English column first for both rows, then the
# Russian/Kyrgyz columns) matches how the shared-string table ends up laid
# out after Excel recomputes it on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C18").Value = "By age (in years) "
$ws.Range("C28").Value = "By education"
$ws.Range("B18").Value = "По возрасту (в годах)"
$ws.Range("A18").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("A28").Value = "Билими боюнча"
$ws.Range("B28").Value = "По образованию"
